$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 115.75
$ws.Range("I12").Value = 112.63636
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 112.63636
$ws.Range("L12").Value = 150
$ws.Range("M12").Value = 57.36364
$ws.Range("N12").Value = -490
$ws.Range("H106").Value = 63494576
$ws.Range("I106").Value = 25643546
$ws.Range("J106").Value = 125002500
$ws.Range("K106").Value = 25643546
$ws.Range("L106").Value = 125002500
$ws.Range("M106").Value = -25642915
$ws.Range("N106").Value = -125003762
$ws.Range("H116").Value = 8441.823
$ws.Range("I116").Value = 12889.444
$ws.Range("J116").Value = 3438.25
$ws.Range("K116").Value = 12889.444
$ws.Range("L116").Value = 3438.25
$ws.Range("M116").Value = -9447.444
$ws.Range("N116").Value = -10322.25
$ws.Range("H132").Value = 7247527.5
$ws.Range("I132").Value = 892.1316
$ws.Range("J132").Value = 41669044
$ws.Range("K132").Value = 2676.3948
$ws.Range("L132").Value = 125007132
$ws.Range("M132").Value = -146.3948
$ws.Range("N132").Value = -125012192
$ws.Range("H135").Value = 1395.9615
$ws.Range("I135").Value = 1098.9778
$ws.Range("K135").Value = 9890.8002
$ws.Range("M135").Value = -7355.8002
$ws.Range("H138").Value = 2168.9158
$ws.Range("I138").Value = 762.92
$ws.Range("J138").Value = 4299.212
$ws.Range("K138").Value = 2288.76
$ws.Range("L138").Value = 12897.636
$ws.Range("M138").Value = 2851.24
$ws.Range("N138").Value = -23177.636
$ws.Range("H141").Value = 1915.1072
$ws.Range("I141").Value = 1439.6818
$ws.Range("K141").Value = 4319.0454
$ws.Range("M141").Value = 860.9546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4832.9478
$ws.Range("I32").Value = 3735.0657
$ws.Range("J32").Value = 9004.9
$ws.Range("K32").Value = 3735.0657
$ws.Range("L32").Value = 9004.9
$ws.Range("M32").Value = -3448.0657
$ws.Range("N32").Value = -9578.9
$ws.Range("H74").Value = 13890000
$ws.Range("I74").Value = 1041.3334
$ws.Range("J74").Value = 55556876
$ws.Range("K74").Value = 1041.3334
$ws.Range("L74").Value = 55556876
$ws.Range("M74").Value = -167.3334
$ws.Range("N74").Value = -55558624
$ws.Range("H77").Value = 13890000
$ws.Range("I77").Value = 1041.3334
$ws.Range("J77").Value = 55556876
$ws.Range("K77").Value = 5206.666999999999
$ws.Range("L77").Value = 277784380
$ws.Range("M77").Value = -838.6669999999995
$ws.Range("N77").Value = -277793116
$ws.Range("H122").Value = 989055.5600000001
$ws.Range("I122").Value = 1605799.8
$ws.Range("J122").Value = 2264.9
$ws.Range("K122").Value = 4817399.4
$ws.Range("L122").Value = 6794.700000000001
$ws.Range("M122").Value = -4814949.4
$ws.Range("N122").Value = -11694.7
$ws.Range("H132").Value = 2229.2334
$ws.Range("I132").Value = 1646.7435
$ws.Range("J132").Value = 3311
$ws.Range("K132").Value = 4940.2305
$ws.Range("L132").Value = 9933
$ws.Range("M132").Value = -2410.2305
$ws.Range("N132").Value = -14993

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4426.9287
$ws.Range("I134").Value = 6142.125
$ws.Range("J134").Value = 2140
$ws.Range("K134").Value = 18426.375
$ws.Range("L134").Value = 6420
$ws.Range("M134").Value = -15891.375
$ws.Range("N134").Value = -11490

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2195.23
$ws.Range("I31").Value = 1405.9464
$ws.Range("J31").Value = 3199.7727
$ws.Range("K31").Value = 1405.9464
$ws.Range("L31").Value = 3199.7727
$ws.Range("M31").Value = -1110.9464
$ws.Range("N31").Value = -3789.7727
$ws.Range("H34").Value = 2195.23
$ws.Range("I34").Value = 1405.9464
$ws.Range("J34").Value = 3199.7727
$ws.Range("K34").Value = 1405.9464
$ws.Range("L34").Value = 3199.7727
$ws.Range("M34").Value = -1203.9464
$ws.Range("N34").Value = -3603.7727
$ws.Range("H58").Value = 2874661.5
$ws.Range("I58").Value = 4274124
$ws.Range("J58").Value = 2080.7896
$ws.Range("K58").Value = 4274124
$ws.Range("L58").Value = 2080.7896
$ws.Range("M58").Value = -4273921
$ws.Range("N58").Value = -2486.7896
$ws.Range("H74").Value = 22164
$ws.Range("J74").Value = 22164
$ws.Range("L74").Value = 22164
$ws.Range("N74").Value = -23912
$ws.Range("H77").Value = 22164
$ws.Range("J77").Value = 22164
$ws.Range("L77").Value = 66492
$ws.Range("N77").Value = -75228
$ws.Range("H132").Value = 4547627
$ws.Range("I132").Value = 6062348
$ws.Range("J132").Value = 3462.3635
$ws.Range("K132").Value = 18187044
$ws.Range("L132").Value = 10387.0905
$ws.Range("M132").Value = -18184514
$ws.Range("N132").Value = -15447.0905
$ws.Range("H134").Value = 10103642
$ws.Range("I134").Value = 13891891
$ws.Range("J134").Value = 1643.7778
$ws.Range("K134").Value = 41675673
$ws.Range("L134").Value = 4931.3334
$ws.Range("M134").Value = -41673138
$ws.Range("N134").Value = -10001.3334
$ws.Range("H136").Value = 2874661.5
$ws.Range("I136").Value = 4274124
$ws.Range("J136").Value = 2080.7896
$ws.Range("K136").Value = 12822372
$ws.Range("L136").Value = 6242.3688
$ws.Range("M136").Value = -12819822
$ws.Range("N136").Value = -11342.3688
$ws.Range("H138").Value = 110000
$ws.Range("J138").Value = 110000
$ws.Range("L138").Value = 110000
$ws.Range("N138").Value = -120280
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 27003.031
$ws.Range("J140").Value = 27003.031
$ws.Range("L140").Value = 27003.031
$ws.Range("N140").Value = -37363.031
$ws.Range("H141").Value = 32379.777
$ws.Range("J141").Value = 32379.777
$ws.Range("L141").Value = 32379.777
$ws.Range("N141").Value = -42739.777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2364253.5
$ws.Range("I5").Value = 208.04
$ws.Range("K5").Value = 624.12
$ws.Range("M5").Value = -512.12
$ws.Range("H26").Value = 50000356
$ws.Range("I26").Value = 103
$ws.Range("J26").Value = 100000610
$ws.Range("K26").Value = 309
$ws.Range("L26").Value = 300001830
$ws.Range("M26").Value = -21
$ws.Range("N26").Value = -300002406
$ws.Range("H74").Value = 20000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 20000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H113").Value = 5455168
$ws.Range("I113").Value = 6250669
$ws.Range("J113").Value = 3333833.2
$ws.Range("K113").Value = 18752007
$ws.Range("L113").Value = 10001499.6
$ws.Range("M113").Value = -18749837
$ws.Range("N113").Value = -10005839.6
$ws.Range("H122").Value = 7385.1
$ws.Range("I122").Value = 394.0909
$ws.Range("J122").Value = 11432.526
$ws.Range("K122").Value = 3546.8181
$ws.Range("L122").Value = 102892.734
$ws.Range("M122").Value = -1096.8181
$ws.Range("N122").Value = -107792.734
$ws.Range("H135").Value = 2364253.5
$ws.Range("I135").Value = 208.04
$ws.Range("K135").Value = 1872.36
$ws.Range("M135").Value = 662.6400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 6672280.5
$ws.Range("I20").Value = 20000000
$ws.Range("J20").Value = 8421
$ws.Range("K20").Value = 20000000
$ws.Range("L20").Value = 8421
$ws.Range("M20").Value = -19999755
$ws.Range("N20").Value = -8911

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10670511
$ws.Range("I132").Value = 13679334
$ws.Range("K132").Value = 41038002
$ws.Range("M132").Value = -41035472

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H62").Value = 3500
$ws.Range("I62").Value = 3500
$ws.Range("K62").Value = 3500
$ws.Range("M62").Value = -2876
$ws.Range("H65").Value = 3500
$ws.Range("I65").Value = 3500
$ws.Range("K65").Value = 17500
$ws.Range("M65").Value = -14380
$ws.Range("H136").Value = 13074541
$ws.Range("I136").Value = 2848.1943
$ws.Range("J136").Value = 44446604
$ws.Range("K136").Value = 8544.582900000001
$ws.Range("L136").Value = 133339812
$ws.Range("M136").Value = -5994.582900000001
$ws.Range("N136").Value = -133344912
